{"js": "// Apply the tracked set of edits to the Ukrainian \"CrisisText Video Scripts\"\n// document:\n//   1. Remove the leftover empty Google-Docs revision-marker content control\n//      (tag \"goog_rdk_0\") that precedes the \"{Make a Routine for Time\n//      Together}\" heading run.\n//   2. Retitle the \"Uncertain Situations\" lesson to \"Crisis\" (heading,\n//      intro sentence and sub-heading), replacing the Ukrainian text with\n//      the new English copy.\n//   3. Drop the trailing \"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c...\" paragraph (and the blank\n//      paragraph right before it) from the \"talk to someone you trust\"\n//      table cell.\n//   4. Clear out the lone reviewer comment (\"\u0422\u0430 \u0441\u0430\u043c\u0430 \u043f\u043e\u0440\u0430\u0434\u0430, \u0449\u043e \u0439 \u0432 AF\").\n\nconst body = context.document.body;\n\n// --- 1. Remove the empty goog_rdk_0 content control -----------------------\nconst controls = context.document.getContentControls();\ncontrols.load(\"items/tag\");\nawait context.sync();\nfor (let i = 0; i < controls.items.length; i++) {\n  if (controls.items[i].tag === \"goog_rdk_0\") {\n    controls.items[i].delete(false);\n  }\n}\nawait context.sync();\n\n// --- 2. Swap the \"Uncertain Situations\" copy for the \"Crisis\" copy --------\nconst headingHits = body.search(\"{Help Children Cope with Uncertain Situations}\", { matchCase: true });\nheadingHits.load(\"items\");\nawait context.sync();\nheadingHits.items.forEach((r) => r.insertText(\"{Help Children Cope with Crisis}\", \"Replace\"));\nawait context.sync();\n\nconst introHits = body.search(\n  \"\u0421\u044c\u043e\u0433\u043e\u0434\u043d\u0456\u0448\u043d\u0456\u0439 \u0443\u0440\u043e\u043a \u043f\u0440\u043e \u0442\u0435, \u044f\u043a \u0434\u043e\u043f\u043e\u043c\u043e\u0433\u0442\u0438 \u0434\u0456\u0442\u044f\u043c \u0443\u043f\u043e\u0440\u0430\u0442\u0438\u0441\u044f \u0437 \u043d\u0435\u0432\u0438\u0437\u043d\u0430\u0447\u0435\u043d\u0456\u0441\u0442\u044e. \",\n  { matchCase: true }\n);\nintroHits.load(\"items\");\nawait context.sync();\nintroHits.items.forEach((r) =>\n  r.insertText(\"Today\u2019s lesson is about helping your children cope with crisis. \", \"Replace\")\n);\nawait context.sync();\n\nconst subheadHits = body.search(\"\u042f\u043a \u0434\u043e\u043f\u043e\u043c\u043e\u0433\u0442\u0438 \u0434\u0456\u0442\u044f\u043c \u0443\u043f\u043e\u0440\u0430\u0442\u0438\u0441\u044f \u0437 \u043d\u0435\u0432\u0438\u0437\u043d\u0430\u0447\u0435\u043d\u0456\u0441\u0442\u044e\", { matchCase: true });\nsubheadHits.load(\"items\");\nawait context.sync();\nsubheadHits.items.forEach((r) => r.insertText(\"Help Children Cope with Crisis\", \"Replace\"));\nawait context.sync();\n\n// --- 3. Delete the \"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c...\" paragraph + the blank one above --\nconst noteHits = body.search(\n  \"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c, \u0449\u043e \u0432\u043e\u043d\u0438 \u0437\u0430\u0432\u0436\u0434\u0438 \u043c\u043e\u0436\u0443\u0442\u044c \u043f\u043e\u0433\u043e\u0432\u043e\u0440\u0438\u0442\u0438 \u0437 \u0432\u0430\u043c\u0438 \u0430\u0431\u043e \u0434\u043e\u0440\u043e\u0441\u043b\u0438\u043c\u0438, \u044f\u043a\u0438\u043c \u0434\u043e\u0432\u0456\u0440\u044f\u044e\u0442\u044c, \u043d\u0430\u0432\u0456\u0442\u044c \u044f\u043a\u0449\u043e \u0457\u043c \u043d\u0435\u043a\u043e\u043c\u0444\u043e\u0440\u0442\u043d\u043e.\",\n  { matchCase: true }\n);\nnoteHits.load(\"items\");\nawait context.sync();\nif (noteHits.items.length > 0) {\n  const targetPara = noteHits.items[0].paragraphs.getFirst();\n  const prevPara = targetPara.getPrevious();\n  targetPara.delete();\n  prevPara.delete();\n  await context.sync();\n}\n\n// --- 4. Clear the reviewer comment -----------------------------------------\nconst comments = context.document.comments;\ncomments.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  comments.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Apply the tracked set of edits to the Ukrainian \"CrisisText Video Scripts\"\n# document:\n#   1. Remove the leftover empty Google-Docs revision-marker content control\n#      (tag \"goog_rdk_0\") that precedes the \"{Make a Routine for Time\n#      Together}\" heading run.\n#   2. Retitle the \"Uncertain Situations\" lesson to \"Crisis\" (heading,\n#      intro sentence and sub-heading), replacing the Ukrainian text with\n#      the new English copy.\n#   3. Drop the trailing \"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c...\" paragraph (and the blank\n#      paragraph right before it) from the \"talk to someone you trust\"\n#      table cell.\n#   4. Clear out the lone reviewer comment (\"\u0422\u0430 \u0441\u0430\u043c\u0430 \u043f\u043e\u0440\u0430\u0434\u0430, \u0449\u043e \u0439 \u0432 AF\"),\n#      including the now-unused comment paragraph style it referenced.\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the empty goog_rdk_0 content control -----------------------\nforeach ($cc in $d.ContentControls) {\n  if ($cc.Tag -eq \"goog_rdk_0\") {\n    $cc.Delete()\n  }\n}\n\n# --- 2. Swap the \"Uncertain Situations\" copy for the \"Crisis\" copy --------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"{Help Children Cope with Uncertain Situations}\"\n$find.Replacement.Text = \"{Help Children Cope with Crisis}\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find.Text = \"\u0421\u044c\u043e\u0433\u043e\u0434\u043d\u0456\u0448\u043d\u0456\u0439 \u0443\u0440\u043e\u043a \u043f\u0440\u043e \u0442\u0435, \u044f\u043a \u0434\u043e\u043f\u043e\u043c\u043e\u0433\u0442\u0438 \u0434\u0456\u0442\u044f\u043c \u0443\u043f\u043e\u0440\u0430\u0442\u0438\u0441\u044f \u0437 \u043d\u0435\u0432\u0438\u0437\u043d\u0430\u0447\u0435\u043d\u0456\u0441\u0442\u044e. \"\n$find.Replacement.Text = \"Today\u2019s lesson is about helping your children cope with crisis. \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find.Text = \"\u042f\u043a \u0434\u043e\u043f\u043e\u043c\u043e\u0433\u0442\u0438 \u0434\u0456\u0442\u044f\u043c \u0443\u043f\u043e\u0440\u0430\u0442\u0438\u0441\u044f \u0437 \u043d\u0435\u0432\u0438\u0437\u043d\u0430\u0447\u0435\u043d\u0456\u0441\u0442\u044e\"\n$find.Replacement.Text = \"Help Children Cope with Crisis\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- 3. Delete the \"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c...\" paragraph + the blank one above --\n$target = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text.Contains(\"\u041d\u0430\u0433\u043e\u043b\u043e\u0441\u0456\u0442\u044c \u0434\u0456\u0442\u044f\u043c, \u0449\u043e \u0432\u043e\u043d\u0438 \u0437\u0430\u0432\u0436\u0434\u0438 \u043c\u043e\u0436\u0443\u0442\u044c \u043f\u043e\u0433\u043e\u0432\u043e\u0440\u0438\u0442\u0438 \u0437 \u0432\u0430\u043c\u0438 \u0430\u0431\u043e \u0434\u043e\u0440\u043e\u0441\u043b\u0438\u043c\u0438, \u044f\u043a\u0438\u043c \u0434\u043e\u0432\u0456\u0440\u044f\u044e\u0442\u044c, \u043d\u0430\u0432\u0456\u0442\u044c \u044f\u043a\u0449\u043e \u0457\u043c \u043d\u0435\u043a\u043e\u043c\u0444\u043e\u0440\u0442\u043d\u043e.\")) {\n    $target = $p\n    break\n  }\n}\nif ($target -ne $null) {\n  $prev = $target.Previous()\n  $target.Range.Delete()\n  $prev.Range.Delete()\n}\n\n# --- 4. Clear the reviewer comment and its leftover style ------------------\nwhile ($d.Comments.Count -gt 0) {\n  $d.Comments.Item(1).Delete()\n}\n\ntry {\n  $style = $d.Styles.Item(\"P68B1DB1-Normal4\")\n  if ($style -ne $null) {\n    $style.Delete()\n  }\n} catch {\n}\n"}
